$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Row 3 on every sheet corresponds to the c3256c5e-... file, whose status
# moves from "Ready for handoff" to "Handback transform failed".
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Record the handback transform error detail for that row on each locale sheet.
$zhcn.Range("L3").Value = "Handback file name: st3iu5g0.yft is different with handoff file name: c3256c5e-b17e-4548-a8a6-bb40d6a77cad.7cc54bf09a31ce50eef2dd8a0b97c9c0280d2861.zh-cn."
$dede.Range("L3").Value = "Handback file name: st3iu5g0.yft is different with handoff file name: c3256c5e-b17e-4548-a8a6-bb40d6a77cad.7cc54bf09a31ce50eef2dd8a0b97c9c0280d2861.de-de."
